$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated prices in column D are plain decimal numbers (e.g. "562.34").
# The sheet stores prices as text, so force those cells to Text format first -
# otherwise Excel would silently convert the assigned string into a number.
$textCells = "D5", "D8", "D10", "D12", "D14", "D15", "D19", "D21", "D22", "D24", "D27", "D31", "D32", "D34", "D36", "D37", "D38", "D40", "D42", "D45", "D47", "D48", "D49"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin stats (price + 1h volume change). Rows 36/37 also swap
# their coin name/link/price/change (RenderToken and FirstDigitalUSD traded places).
$ws.Range("D2").Value = '62.263.81'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '2.420.67'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '562.34'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("E6").Value = '  +1.45%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  +1.24%  '
$ws.Range("D9").Value = '2.417.61'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").Value = '0.109'
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("E11").Value = '  -2.07%  '
$ws.Range("D12").Value = '5.34'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '25.83'
$ws.Range("E14").Value = '  +0.24%  '
$ws.Range("D15").Value = '0.0000176'
$ws.Range("E15").Value = '  +0.09%  '
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("D17").Value = '62.127.70'
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("D18").Value = '2.418.14'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '11.29'
$ws.Range("E19").Value = '  +2.32%  '
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '324.25'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '6.82'
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '65.82'
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("E25").Value = '  -3.96%  '
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").Value = '579.87'
$ws.Range("E27").Value = '  +7.72%  '
$ws.Range("D28").Value = '0.0₃0953'
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("D29").Value = '2.538.49'
$ws.Range("E29").Value = '  +1.47%  '
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("D31").Value = '8.27'
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("D32").Value = '1.45'
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("D34").Value = '1.89'
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D37").Value = '5.66'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("D40").Value = '152.10'
$ws.Range("E40").Value = '  +3.53%  '
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("D42").Value = '1.82'
$ws.Range("E42").Value = '  -5.61%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("E44").Value = '  +4.65%  '
$ws.Range("D45").Value = '149.79'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("D47").Value = '0.0537'
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("D48").Value = '20.21'
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("D49").Value = '0.595'
$ws.Range("E49").Value = '  +1.60%  '
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("E51").Value = '  +1.42%  '

# Restore the default "Normal" style on the text-forced cells so only their
# value changed - the Text number-format was only needed during assignment.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
